# Update wording for two closed-bug titles in the "closed bugs in last
# iteration" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("closed bugs in last iteration")

# Bug 2929769 (row 5): reworded title
$ws.Range("B5").Value = "Loading mask artifacat when deleting packages in Packages grid from Feed details"

# Bug 2901954 (row 12): reworded title
$ws.Range("B12").Value = "Update Grafana dependencies to resolve public CVEs"

# Keep the active cell selection consistent with the edited file
$ws.Range("B17").Select()
